$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the current column H (salary) so that
# H/I/J (salary, date updated, date created) shift right to I/J/K.
$ws.Columns.Item(8).Insert()

# New header for inserted column H: "Vice-fixed"
$ws.Range("H1").Value = "Vice-fixed"

# Fill the new column with a TRIM() cleanup formula mirroring the existing
# PROPER() "President-fixed" column D.
$ws.Range("H2:H47").Formula = "=TRIM(G2)"

# Header fill: medium blue (matches the "salary" style header look)
$ws.Range("H1").Interior.ThemeColor = 5
$ws.Range("H1").Interior.TintAndShade = -0.4

# Data column fill: light blue
$ws.Range("H2:H47").Interior.ThemeColor = 5
$ws.Range("H2:H47").Interior.TintAndShade = 0.8

# Party column (F) gets new fills too per the diff.
$ws.Range("F1").Interior.ThemeColor = 1
$ws.Range("F1").Interior.TintAndShade = -0.25

$ws.Range("F2:F47").Interior.ThemeColor = 3
$ws.Range("F2:F47").Interior.TintAndShade = 0.6

# Update the view so the top-left visible cell is B1 and the active
# selection is H1 (matches the diff's sheetView change).
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("H1").Select()
